# Updates Leve profit-calculation figures (currentAveragePrice / LevePrice /
# LeveProfit columns H-N) across the per-job sheets. Values refresh from the
# latest market-board pull; a handful of rows also gain/lose the Profit
# columns (M/N) when the underlying price data goes from present to missing
# or vice versa.
#
# $null in a column's value list means "clear that cell" (ClearContents),
# matching rows where the source data no longer yields a LeveProfit figure.

$wb = $excel.ActiveWorkbook

$sheetEdits = @{
    "ALC" = @{
        2   = @{ H=484.2;       I=334.7143;          J=833;         K=334.7143;          L=833;         M=-221.7143;          N=-1059 }
        17  = @{ H=4180.2812;                        J=4180.2812;                  L=12540.8436;                               N=-12876.8436 }
        132 = @{ H=736.0267;    I=650.92957;          J=2246.5;      K=1952.78871;        L=6739.5;      M=577.21129;          N=-11799.5 }
        135 = @{ H=35715610;    I=13890080;           J=166668770;   K=125010720;         L=1500018930;  M=-125008185;         N=-1500024000 }
    }
    "ARM" = @{
        32  = @{ H=10269.837;   I=7900.3047;          J=29700;       K=7900.3047;         L=29700;       M=-7613.3047;         N=-30274 }
        74  = @{ H=97488.37;    I=143421.2;           J=9156;        K=143421.2;          L=9156;        M=-142547.2;          N=-10904 }
        77  = @{ H=97488.37;    I=143421.2;           J=9156;        K=717106;            L=45780;       M=-712738;            N=-54516 }
        122 = @{ H=8930714;     I=2285.7144;          J=17859142;    K=6857.1432;         L=53577426;    M=-4407.1432;         N=-53582326 }
    }
    "BSM" = @{
        134 = @{ H=21799.79;                          J=104154.2;                   L=312462.6;                                N=-317532.6 }
    }
    "CRP" = @{
        31  = @{ H=17998.8;     I=0;                  J=17998.8;     K=0;                 L=17998.8;     M=$null;              N=-18588.8 }
        34  = @{ H=17998.8;     I=0;                  J=17998.8;     K=0;                 L=17998.8;     M=$null;              N=-18402.8 }
        107 = @{ H=1364.9333;   I=1631.6666;          J=964.8333;    K=1631.6666;         L=964.8333;    M=288.3334;           N=-4804.8333 }
        109 = @{ H=0;           I=0;                  J=0;           K=0;                 L=0;           M=$null;              N=$null }
    }
    "CUL" = @{
        4   = @{ H=1450;        I=716.6667;                          K=2150.0001;                        M=-2038.0001 }
        35  = @{ H=1400;        I=1000;                              K=3000;                              M=-2712 }
        122 = @{ H=993.85187;   I=555.5333000000001;  J=1541.75;     K=4999.7997;         L=13875.75;    M=-2549.7997;         N=-18775.75 }
        130 = @{ H=5147.2;      I=1620;                J=7498.6665;   K=4860;              L=22495.9995;  M=160;                N=-32535.9995 }
        133 = @{ H=5115.2354;   I=5734.75;             J=4924.615;    K=17204.25;          L=14773.845;   M=-12144.25;          N=-24893.845 }
        137 = @{ H=45459356;    I=55558428;                          K=166675284;                        M=-166670184 }
    }
    "GSM" = @{
        14  = @{ H=18000000;    I=18000000;                          K=18000000;                         M=-17999832 }
        19  = @{ H=0;                                  J=0;                        L=0;                                       N=$null }
        102 = @{ H=5887.1904;   I=5440.3076;           J=6613.375;    K=5440.3076;         L=6613.375;    M=-3818.3076;         N=-9857.375 }
        103 = @{ H=25000;                              J=25000;                    L=25000;                                   N=-27344 }
    }
    "LTW" = @{
        22  = @{ H=1974.75;     I=1900;                J=1999.6666;   K=1900;              L=1999.6666;   M=-1605;              N=-2589.6666 }
        27  = @{ H=1974.75;     I=1900;                J=1999.6666;   K=1900;              L=1999.6666;   M=-1793;              N=-2213.6666 }
        46  = @{ H=1466.6666;   I=1200;                J=2000;        K=1200;              L=2000;        M=-1012;              N=-2376 }
        93  = @{ H=6040;        I=5300;                J=9000;        K=5300;              L=9000;        M=-4052;              N=-11496 }
        132 = @{ H=7041;        I=8771.200000000001;                 K=26313.6;                          M=-23783.6 }
        136 = @{ H=5963.718;    I=4456.75;             J=7550;        K=13370.25;          L=22650;       M=-10820.25;          N=-27750 }
    }
    "WVR" = @{
        122 = @{ H=8876.444;    I=3496.3333;           J=19636.666;   K=10488.9999;        L=58909.99800000001; M=-8038.999899999999; N=-63809.99800000001 }
        136 = @{ H=4720.3335;   I=2203.8;              J=7561.5806;   K=6611.400000000001; L=22684.7418;  M=-4061.400000000001; N=-27784.7418 }
        141 = @{ H=42967;                              J=44335.555;                L=44335.555;                                N=-54695.555 }
    }
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetEdits[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $cols = $rows[$rowNum]
        foreach ($colLetter in $cols.Keys) {
            $cellValue = $cols[$colLetter]
            $addr = "$colLetter$rowNum"
            if ($null -eq $cellValue) {
                $ws.Range($addr).ClearContents()
            } else {
                $ws.Range($addr).Value = $cellValue
            }
        }
    }
}

Write-Host "Updated Leve profit figures across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."
